$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.232.15'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.68%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.536.33'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.83%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '543.74'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.95'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.43%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.572'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.564.83'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.83%  '
$ws.Range("E10").Value = '  +0.54%  '
$ws.Range("E11").Value = '  +0.88%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.53'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.75%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.363'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.990.23'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.80'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -2.92%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '60.024.93'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.47%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000143'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +2.50%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.552.15'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.60%  '
$ws.Range("E19").Value = '  -2.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.32'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '327.71'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.40%  '
$ws.Range("E22").Value = '  -0.69%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.93'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +2.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '62.59'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.442'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.31%  '
$ws.Range("E26").Value = '  +2.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.992'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.01'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.77%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.06'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0797'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.64%  '
$ws.Range("E31").Value = '  -0.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.23'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.63%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '162.48'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.87%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.47'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +2.84%  '
$ws.Range("E35").Value = '  +0.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.77'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.46'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.63'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.71'
$ws.Range("D39").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.14'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '302.35'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -4.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.838'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.74'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.11%  '
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.995'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.609'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.84'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.37%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '19.05'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.09%  '
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0938'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.41%  '
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0938'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.46%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0520'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.78%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0229'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.29%  '
